$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the Timestamp value on the existing last row (row 6) ---
# Tiny floating-point correction picked up by the "get" API vs the old "post" API.
$ws.Range("A6").Value = 45729.51109761574

# --- Shared recommendation text reused by the three new rows ---
$calciumText = "`nRecommendations for calcium Deficiency:`nTofu, raw, firm, prepared with calcium sulfate`nCheese, Mexican, blend, reduced fat`nCheese, cheddar, nonfat or fat free`nCheese, swiss, low fat`nCheese, swiss, low sodium`nCheese, mozzarella, part skim milk`nCheese, gruyere`nCheese, monterey`nCheese, port de salut`nCheese, swiss`nCheese, swiss`nCheese, provolone, sliced`nCheese, provolone, reduced fat`nCheese, monterey jack, solid`nCheese, low-sodium, cheddar or colby`nCheese, muenster`nCheese, mozzarella, low sodium`nCheese, provolone`nCheese, monterey, low fat`nCheese, brick`nCheese, mexican, queso asadero`nCheese, colby`nCheese, Mexican blend`nCheese, Swiss, nonfat or fat free`nCheese, queso fresco, solid`nCheese, cheddar`nCheese, mexican, queso chihuahua`nCheese, cheddar, sharp, sliced`nCheese, cheddar`nCheese, white, queso blanco`nCheese, mozzarella, nonfat`nCheese, cheddar, reduced fat`nCheese, tilsit`nCheese, parmesan, grated, refrigerated`nCheese, cheshire`nCheese, parmesan, hard`nCheese, caraway`nImitation cheese, american or cheddar, low cholesterol`nCheese, fontina`nCheese, mexican, queso anejo"

$bmiCategory = "Normal weight - Maintain a balanced diet and exercise."

# --- New rows 7, 8, 9: three "get" API calls logged right after row 6 ---
$timestamps = @(45730.43464201389, 45730.43493046296, 45730.43550134636)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $r = 7 + $i

    # Match the Timestamp column's datetime number format used by the rows above.
    $ws.Range("A$r").NumberFormat = $ws.Range("A6").NumberFormat
    $ws.Range("A$r").Value = $timestamps[$i]

    $ws.Range("B$r").Value = "Aarti"
    $ws.Range("C$r").Value = 25
    $ws.Range("D$r").Value = "Female"
    $ws.Range("E$r").Value = 50
    $ws.Range("F$r").Value = 1.5
    $ws.Range("G$r").Value = 22.22
    $ws.Range("H$r").Value = $bmiCategory
    $ws.Range("I$r").Value = "Veg"
    $ws.Range("J$r").Value = "calcium"
    $ws.Range("K$r").Value = $calciumText

    # Writing the long multi-line recommendation auto-expands the row height;
    # AutoFit it back so no explicit/custom row height is persisted.
    $ws.Rows.Item($r).AutoFit()
}
